$d = $word.ActiveDocument

# --- Fill in the two title/subtitle content controls -----------------
# Setting Range.Text on a content control that is still showing its
# placeholder text replaces the placeholder runs with a single run of
# real text and clears the w:showingPlcHdr flag automatically.
for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Tag -eq "ccDocumentTitle") {
        $cc.Range.Text = "asdfsf"
    } elseif ($cc.Tag -eq "ccDocumentSubtitle") {
        $cc.Range.Text = "sadfasdfasdf"
    }
}

# --- Remove the stray _GoBack bookmark --------------------------------
# Word drops this automatically managed bookmark when it is no longer
# needed; removing it also causes the remaining bookmark ids to be
# renumbered (0-based) on save.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
